$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match row (row 95) appended right after the last existing data row (row 94).
# Copy formatting (styles) from row 94 down to row 95 first so the new row
# inherits the same per-column styles (bold/bordered index column, date-time
# number format column, etc.) without minting any new style entries.
$ws.Range("A94:V94").Copy()
$ws.Range("A95:V95").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Plain text / numeric cells.
$ws.Cells.Item(95,1).Value = 94
$ws.Cells.Item(95,2).Value = "ecuador"
$ws.Cells.Item(95,3).Value = "liga-pro"

# Column D holds "2023" as TEXT (matching every other row), but a plain
# Value assignment of a purely-numeric-looking string gets auto-coerced to a
# number. Route it through a TEXT() formula and paste back as a value so the
# cell lands as a literal string without adding a quotePrefix/number-format
# style to the workbook.
$ws.Cells.Item(95,4).Formula = '=TEXT(2023,"0")'
$ws.Range("D95").Copy()
$ws.Range("D95").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$ws.Cells.Item(95,5).Value = 45236
$ws.Cells.Item(95,6).Value = "Barcelona SC"
$ws.Cells.Item(95,7).Value = 3
$ws.Cells.Item(95,8).Value = "EL Nacional"
$ws.Cells.Item(95,9).Value = 2
$ws.Cells.Item(95,10).Value = 1.79
$ws.Cells.Item(95,11).Value = "31/10/2023 01:12"
$ws.Cells.Item(95,12).Value = 1.71
$ws.Cells.Item(95,13).Value = "05/11/2023 23:56"
$ws.Cells.Item(95,14).Value = 4.15
$ws.Cells.Item(95,15).Value = "31/10/2023 01:12"
$ws.Cells.Item(95,16).Value = 4.06
$ws.Cells.Item(95,17).Value = "05/11/2023 23:56"
$ws.Cells.Item(95,18).Value = 3.6
$ws.Cells.Item(95,19).Value = "31/10/2023 01:12"
$ws.Cells.Item(95,20).Value = 4.61
$ws.Cells.Item(95,21).Value = "05/11/2023 23:56"
$ws.Cells.Item(95,22).Value = "https://www.betexplorer.com/football/ecuador/liga-pro/barcelona-sc-el-nacional/t2KyI6KN/"
